# Fixed inability to pass in duplicate asset descriptions as input
#
# The "partial overhead" formulas in column B (*_partial_oh) and column G
# (budgeted_partial_oh) were incorrectly subtracting the labour/contracts/parts
# breakdown columns a second time (e.g. "=C2-D2-E2-F2" instead of simply
# "=C2-D2"), which broke whenever a duplicate asset description/line item was
# supplied. Simplify those formulas on both data sheets.

$wb = $excel.ActiveWorkbook

$wsW = $wb.Worksheets.Item("BME_DI_BCCW")
$wsC = $wb.Worksheets.Item("BME_DI_BCC")

foreach ($ws in @($wsW, $wsC)) {
    for ($r = 2; $r -le 6; $r++) {
        $ws.Range("B$r").Formula = "=C$r-D$r"
        $ws.Range("G$r").Formula = "=H$r-I$r"
    }
}

# --- Update selections / active sheet to match the saved UI state ---
# BME_DI_BCC is no longer the active tab, and its cursor ends up on G4.
$wsC.Activate()
$wsC.Range("G4").Select()

# BME_DI_BCCW becomes the active (selected) tab, with the cursor on G4.
$wsW.Activate()
$wsW.Range("G4").Select()
